# Remove the post entry "「結婚前、結婚後」" which occupied row 804.
# Deleting the entire row shifts every subsequent row (805-829) up by one
# (805->804, ..., 829->828), matching the published diff, and Excel
# automatically shrinks the sheet's used-range dimension from C829 to C828.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(804).EntireRow.Delete()
